# Generate Report for Handoff
# Replaces the two previously-handed-back localization rows with a single
# fresh "ready for handoff" entry (new source file 786ce09c-...-md, new
# target xliff hashes) and a second row (ffff3837dc51-...-md) that shares
# the same new handoff package, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

function Remove-HyperlinkAt($sheet, $addr) {
    foreach ($h in $sheet.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.Delete()
            return
        }
    }
}

function Set-HyperlinkDisplay($sheet, $addr, $text) {
    foreach ($h in $sheet.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $text
            return
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.md"
$ov.Range("B2").Value = "e2e\786ce09c-c2f1-4abb-bc18-5632705c6788.md"
Set-HyperlinkDisplay $ov '$B$2' "e2e\786ce09c-c2f1-4abb-bc18-5632705c6788.md"
$ov.Range("C2").Value = ".md"
$ov.Range("D2").Value = ""
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-08-21 07:07:52"

$ov.Range("A3").Value = "ffff3837dc51-0456-44ef-a254-0496531384a2.md"
$ov.Range("B3").Value = "e2e\ffff3837dc51-0456-44ef-a254-0496531384a2.md"
Set-HyperlinkDisplay $ov '$B$3' "e2e\ffff3837dc51-0456-44ef-a254-0496531384a2.md"
$ov.Range("C3").Value = ".md"
$ov.Range("D3").Value = ""
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-21 07:07:52"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

Remove-HyperlinkAt $zh '$I$2'
Remove-HyperlinkAt $zh '$I$3'

$zh.Range("A2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.md"
Set-HyperlinkDisplay $zh '$A$2' "786ce09c-c2f1-4abb-bc18-5632705c6788.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("D2").Value = "e2e"
$zh.Range("E2").Value = "ht"
$zh.Range("F2").Value = "'False"
$zh.Range("G2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.eb3a7322583a465fbc008f875572b03f182d6465.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-21 07:07:48"
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"

$zh.Range("A3").Value = "ffff3837dc51-0456-44ef-a254-0496531384a2.md"
Set-HyperlinkDisplay $zh '$A$3' "ffff3837dc51-0456-44ef-a254-0496531384a2.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "ht"
$zh.Range("F3").Value = "'True"
$zh.Range("G3").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.eb3a7322583a465fbc008f875572b03f182d6465.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-21 07:07:48"
$zh.Range("I3").Value = ""
$zh.Range("J3").Value = ""
$zh.Range("K3").Value = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

Remove-HyperlinkAt $de '$I$2'
Remove-HyperlinkAt $de '$I$3'

$de.Range("A2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.md"
Set-HyperlinkDisplay $de '$A$2' "786ce09c-c2f1-4abb-bc18-5632705c6788.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("D2").Value = "e2e"
$de.Range("E2").Value = "ht"
$de.Range("F2").Value = "'False"
$de.Range("G2").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.eb3a7322583a465fbc008f875572b03f182d6465.de-de.xlf"
$de.Range("H2").Value = "2016-08-21 07:07:52"
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"

$de.Range("A3").Value = "ffff3837dc51-0456-44ef-a254-0496531384a2.md"
Set-HyperlinkDisplay $de '$A$3' "ffff3837dc51-0456-44ef-a254-0496531384a2.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "ht"
$de.Range("F3").Value = "'True"
$de.Range("G3").Value = "786ce09c-c2f1-4abb-bc18-5632705c6788.eb3a7322583a465fbc008f875572b03f182d6465.de-de.xlf"
$de.Range("H3").Value = "2016-08-21 07:07:52"
$de.Range("I3").Value = ""
$de.Range("J3").Value = ""
$de.Range("K3").Value = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Column width autofit (content got shorter after the edits above)
# ---------------------------------------------------------------------
$ov.Range("E:F").Columns.AutoFit()
$zh.Range("C:C").Columns.AutoFit()
$zh.Range("I:J").Columns.AutoFit()
$de.Range("C:C").Columns.AutoFit()
$de.Range("I:J").Columns.AutoFit()

Write-Output "done"
